# Normalize the "Recorded By" (column G) cell values: any entry that is
# exactly "System" (case-sensitive) gets moved to the end of the
# comma-separated list, while the relative order of every other entry
# (including a lowercase "system") is preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $raw = $cell.Value2

    if ($raw -eq $null) {
        continue
    }

    $parts = $raw -split ", "

    if ($parts.Length -le 1) {
        continue
    }

    $nonSystem = @()
    $systemParts = @()

    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $systemParts += $p
        } else {
            $nonSystem += $p
        }
    }

    if ($systemParts.Length -eq 0) {
        continue
    }

    $reordered = $nonSystem + $systemParts
    $newValue = $reordered -join ", "

    if (-not $newValue.Equals($raw)) {
        $cell.Value = $newValue
    }
}

"done"
